$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Third inline picture (the one without w:noProof yet) gets w:noProof
#    added to its run properties.
# ---------------------------------------------------------------------------
$shp = $d.InlineShapes.Item(3)
$shp.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# 2) Q1 answer: the S3 sentence gets split into three runs, and
#    "Can connect" becomes "It can connect".
# ---------------------------------------------------------------------------
$f = $d.Content
$f.Find.Execute(" is a storage facility that can be accessed anywhere, it is storage for the internet. S3 is an object storage system and not a file system. Can connect directly to the internet. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Range($f.Start, $f.End)

$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> is a storage facility that can be accessed anywhere, it is storage for the internet. S3 is an object storage</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> system and not a file system. It c</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">an connect directly to the internet. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Q2 answer: the "I would consider..." sentence is re-typed (creating a
#    run for each edited chunk), "needed either" -> "need either", and
#    " than available on my local machine" is inserted before the final
#    sentence. The _GoBack bookmark now sits right before that final
#    sentence.
# ---------------------------------------------------------------------------
$f = $d.Content
$f.Find.Execute("I would consider a cloud infrastructure for my data science tasks when I needed either more computation power than my machine provided and/or when I needed more storage for my data.  The cloud services can provide additional space for data storage and save time by cutting down on computation times of the job or simply speeding up reading/writing from a database.  I would also consider using an EBS storage to preserve data, version data and not worry about data loss due to machine crash without a backup.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Range($f.Start, $f.End)

$xml3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>I would consider a cloud infras</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>tructure for my data science tas</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">ks when I </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>need</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> either more computation power than my machine provided and/or when I needed more storage for my data</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> than available on my local machine</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>.  The cloud services can provide additional space for data storage and save time by cutting down on computation times of the job or simply speeding up reading/writing from a database.  I would also consider using an EBS storage to preserve data, version data and not worry about data loss due to machine crash without a backup.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 4) Q3 answer: the spot-instance sentence is merged back into a single run
#    (the _GoBack bookmark that used to live here moved to step 3 above).
# ---------------------------------------------------------------------------
$f = $d.Content
$f.Find.Execute("In a spot instance, you bid on an instance. You choose the price you are willing to pay for that configuration.  The price of each instance depends on supply and demand of the configuration.  When the price drops to at or below your specified bid, the instance will turn on and be yours to use until the price rises above your bid again.  These can save money but the application you are using the instance to compute must be able to handle interruptions and should not be time sensitive since it is tough to predict when the instance will be available and for how long. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $d.Range($f.Start, $f.End)

$xml4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve">In a spot instance, you bid on an instance. You choose the price you are willing to pay for that configuration.  The price of each instance depends on supply and demand of the configuration.  When the price drops to at or below your specified bid, the instance will turn on and be yours to use until the price rises above your bid again.  These can save money but the application you are using the instance to compute must be able to handle interruptions and should not be time sensitive since it is tough to predict when the instance will be available and for how long. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml4)

Write-Output "done"
